# Seed dummy "auth user groups" transactions into Sheet1.
# Row 1 (headers: COMPANY ID / GROUP NAME) is left untouched.
# Row 2's company id is refreshed, and rows 3-15 of new dummy records
# are appended below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("EMP-ID-06", "EMPLOYEE ADMIN"),
    @("EMP-ID-25", "PURCHASE REQUEST REQUESTOR"),
    @("EMP-ID-26", "PURCHASE REQUEST REQUESTOR"),
    @("EMP-ID-27", "PURCHASE REQUEST REQUESTOR, purchase request approver"),
    @("EMP-ID-20", "purchase request approver"),
    @("EMP-ID-21", "purchase request approver"),
    @("EMP-ID-22", "purchase request approver"),
    @("EMP-ID-23", "purchase request approver"),
    @("EMP-ID-24", "purchase request approver"),
    @("EMP-ID-28", "SALES REPRESENTATIVE"),
    @("EMP-ID-29", "SALES REPRESENTATIVE"),
    @("EMP-ID-30", "SALES REPRESENTATIVE"),
    @("EMP-ID-31", "SALES REPRESENTATIVE"),
    @("EMP-ID-32", "SALES REPRESENTATIVE, SALES ADMIN")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $row = $row + 1
}

# Matches the author's final selection left in the sheet after entering data.
$ws.Range("B19").Select()
